$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content looks like a plain number need to be forced
# to remain text (matching the source inlineStr cells) by marking the range
# as Text before assignment, then resetting the style so no stray formatting
# is left behind on the cell.

$ws.Range('D2').Value = '75.643.88'
$ws.Range('D3').Value = '2.711.68'
$ws.Range('E3').Value = '  +11.72%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '189.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +13.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '590.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.543'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.86%  '
$ws.Range('E9').Value = '  +17.10%  '
$ws.Range('D10').Value = '2.710.36'
$ws.Range('E10').Value = '  +11.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.163'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.362'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.78'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.183.25'
$ws.Range('E14').Value = '  +10.82%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '75.435.63'
$ws.Range('E15').Value = '  +9.09%  '
$ws.Range('E16').Value = '  +7.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.86'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +12.33%  '
$ws.Range('D18').Value = '2.698.14'
$ws.Range('E18').Value = '  +10.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +34.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +13.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +17.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.13%  '
$ws.Range('E24').Value = '  +4.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.62'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.38%  '
$ws.Range('E29').Value = '  +10.89%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').Value = '0.0₃0980'
$ws.Range('E31').Value = '  +15.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '526.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +16.27%  '
$ws.Range('E33').Value = '  +14.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.92'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.28%  '
$ws.Range('E35').Value = '  +10.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.121'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.21'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.39'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '174.26'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +28.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.09'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +15.58%  '
$ws.Range('E44').Value = '  +13.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.335'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +12.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +16.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '39.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0852'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +18.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.70'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.547'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.88%  '
